$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.680.12"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.807.17"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'225.10"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'0.604"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'40.24"
$ws.Range("E8").Value = "  +10.58%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.0672"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("D12").Value = "2.068.76"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.808.41"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "'10.89"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("D15").Value = "'0.635"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "34.682.34"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "'67.99"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "'241.21"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").Value = "'172.01"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -5.62%  "
$ws.Range("D27").Value = "'17.48"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "'3.78"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "'3.85"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "1.306.20"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'14.85"
$ws.Range("E38").Value = "  +10.71%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0188"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.34"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'84.25"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("D43").Value = "'2.43"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'2.79"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'0.943"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E46").Value = "  +5.27%  "
$ws.Range("D47").Value = "1.968.17"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'101.11"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "'0.0612"
$ws.Range("E51").Value = "  +0.59%  "
